# Update "想去人数" (interested-attendee count) values on the three affected
# sheets: "展览" (Exhibitions), "演出" (Performances), and "全部类型" (All types).
# "本地生活" (Local life) is untouched.

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("F2").Value = 40
$wsExhibitions.Range("F6").Value = 369
$wsExhibitions.Range("F13").Value = 433
$wsExhibitions.Range("F14").Value = 1629
$wsExhibitions.Range("F18").Value = 1392
$wsExhibitions.Range("F19").Value = 269
$wsExhibitions.Range("F21").Value = 1126
$wsExhibitions.Range("F22").Value = 402
$wsExhibitions.Range("F24").Value = 3495
$wsExhibitions.Range("F25").Value = 692
$wsExhibitions.Range("F27").Value = 1547

$wsPerformances = $wb.Worksheets.Item("演出")
$wsPerformances.Range("F8").Value = 22
$wsPerformances.Range("F13").Value = 17

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 40
$wsAllTypes.Range("F13").Value = 22
$wsAllTypes.Range("F16").Value = 369
$wsAllTypes.Range("F23").Value = 433
$wsAllTypes.Range("F24").Value = 1629
$wsAllTypes.Range("F28").Value = 1392
$wsAllTypes.Range("F29").Value = 269
$wsAllTypes.Range("F33").Value = 1126
$wsAllTypes.Range("F34").Value = 402
$wsAllTypes.Range("F36").Value = 3495
$wsAllTypes.Range("F37").Value = 692
$wsAllTypes.Range("F39").Value = 1547
$wsAllTypes.Range("F41").Value = 17
